$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4
$ws.Range("A2").Value = "s1"
$ws.Range("B2").Value = 400

$ws.Range("A3").Value = "s2"
$ws.Range("B3").Value = 300

$ws.Range("A4").Value = "s3"
$ws.Range("B4").Value = 0

# Add new rows 5-10
$ws.Range("A5").Value = "s4"
$ws.Range("B5").Value = 700

$ws.Range("A6").Value = "x1"
$ws.Range("B6").Value = 0

$ws.Range("A7").Value = "x2"
$ws.Range("B7").Value = 2000

$ws.Range("A8").Value = "x3"
$ws.Range("B8").Value = 0

$ws.Range("A9").Value = "x4"
$ws.Range("B9").Value = 0

$ws.Range("A10").Value = "Z (Função Objetivo)"
$ws.Range("B10").Value = 97600
